$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new '2022-Q3' worksheet right after the total summary
#    sheet, i.e. immediately before '2022-Q2'.
# ------------------------------------------------------------------
$anchorSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($anchorSheet, $null)
$q3.Name = "2022-Q3"

# Header row (bold, centered, matches sibling quarter sheets)
$hc = $q3.Cells.Item(1, 2)
$hc.Value = '基金代码'
$hc.Font.Bold = $true
$hc.HorizontalAlignment = -4108
$hc.VerticalAlignment = -4160
$hc = $q3.Cells.Item(1, 3)
$hc.Value = '基金名称'
$hc.Font.Bold = $true
$hc.HorizontalAlignment = -4108
$hc.VerticalAlignment = -4160
$hc = $q3.Cells.Item(1, 4)
$hc.Value = '基金规模'
$hc.Font.Bold = $true
$hc.HorizontalAlignment = -4108
$hc.VerticalAlignment = -4160
$hc = $q3.Cells.Item(1, 5)
$hc.Value = '股票总仓位'
$hc.Font.Bold = $true
$hc.HorizontalAlignment = -4108
$hc.VerticalAlignment = -4160
$hc = $q3.Cells.Item(1, 6)
$hc.Value = '仓位占比'
$hc.Font.Bold = $true
$hc.HorizontalAlignment = -4108
$hc.VerticalAlignment = -4160
$hc = $q3.Cells.Item(1, 7)
$hc.Value = '持有市值(亿元)'
$hc.Font.Bold = $true
$hc.HorizontalAlignment = -4108
$hc.VerticalAlignment = -4160
$hc = $q3.Cells.Item(1, 8)
$hc.Value = '仓位排名'
$hc.Font.Bold = $true
$hc.HorizontalAlignment = -4108
$hc.VerticalAlignment = -4160

# Data rows
$ac = $q3.Cells.Item(2, 1)
$ac.Value = 0
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(2, 2); $c.NumberFormat = "@"; $c.Value = '002510'
$q3.Cells.Item(2, 3).Value = '申万菱信中证500指数增强A'
$c = $q3.Cells.Item(2, 4); $c.NumberFormat = "@"; $c.Value = '3.74'
$c = $q3.Cells.Item(2, 5); $c.NumberFormat = "@"; $c.Value = '93.28'
$c = $q3.Cells.Item(2, 6); $c.NumberFormat = "@"; $c.Value = '1.85'
$c = $q3.Cells.Item(2, 7); $c.NumberFormat = "@"; $c.Value = '0.0692'
$q3.Cells.Item(2, 8).Value = 6
$ac = $q3.Cells.Item(3, 1)
$ac.Value = 1
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(3, 2); $c.NumberFormat = "@"; $c.Value = '009992'
$q3.Cells.Item(3, 3).Value = '景顺长城量化成长演化混合'
$c = $q3.Cells.Item(3, 4); $c.NumberFormat = "@"; $c.Value = '2.30'
$c = $q3.Cells.Item(3, 5); $c.NumberFormat = "@"; $c.Value = '91.84'
$c = $q3.Cells.Item(3, 6); $c.NumberFormat = "@"; $c.Value = '2.44'
$c = $q3.Cells.Item(3, 7); $c.NumberFormat = "@"; $c.Value = '0.0561'
$q3.Cells.Item(3, 8).Value = 10
$ac = $q3.Cells.Item(4, 1)
$ac.Value = 2
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(4, 2); $c.NumberFormat = "@"; $c.Value = '006048'
$q3.Cells.Item(4, 3).Value = '长城中证500指数增强A'
$c = $q3.Cells.Item(4, 4); $c.NumberFormat = "@"; $c.Value = '2.12'
$c = $q3.Cells.Item(4, 5); $c.NumberFormat = "@"; $c.Value = '94.78'
$c = $q3.Cells.Item(4, 6); $c.NumberFormat = "@"; $c.Value = '1.87'
$c = $q3.Cells.Item(4, 7); $c.NumberFormat = "@"; $c.Value = '0.0396'
$q3.Cells.Item(4, 8).Value = 7
$ac = $q3.Cells.Item(5, 1)
$ac.Value = 3
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(5, 2); $c.NumberFormat = "@"; $c.Value = '159804'
$q3.Cells.Item(5, 3).Value = '国寿安保国证创业板中盘精选88ETF'
$c = $q3.Cells.Item(5, 4); $c.NumberFormat = "@"; $c.Value = '1.10'
$c = $q3.Cells.Item(5, 5); $c.NumberFormat = "@"; $c.Value = '98.91'
$c = $q3.Cells.Item(5, 6); $c.NumberFormat = "@"; $c.Value = '1.87'
$c = $q3.Cells.Item(5, 7); $c.NumberFormat = "@"; $c.Value = '0.0206'
$q3.Cells.Item(5, 8).Value = 8
$ac = $q3.Cells.Item(6, 1)
$ac.Value = 4
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(6, 2); $c.NumberFormat = "@"; $c.Value = '007413'
$q3.Cells.Item(6, 3).Value = '长城中证500指数增强C'
$c = $q3.Cells.Item(6, 4); $c.NumberFormat = "@"; $c.Value = '1.06'
$c = $q3.Cells.Item(6, 5); $c.NumberFormat = "@"; $c.Value = '94.78'
$c = $q3.Cells.Item(6, 6); $c.NumberFormat = "@"; $c.Value = '1.87'
$c = $q3.Cells.Item(6, 7); $c.NumberFormat = "@"; $c.Value = '0.0198'
$q3.Cells.Item(6, 8).Value = 7
$ac = $q3.Cells.Item(7, 1)
$ac.Value = 5
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(7, 2); $c.NumberFormat = "@"; $c.Value = '007795'
$q3.Cells.Item(7, 3).Value = '申万菱信中证500指数增强C'
$c = $q3.Cells.Item(7, 4); $c.NumberFormat = "@"; $c.Value = '0.87'
$c = $q3.Cells.Item(7, 5); $c.NumberFormat = "@"; $c.Value = '93.28'
$c = $q3.Cells.Item(7, 6); $c.NumberFormat = "@"; $c.Value = '1.85'
$c = $q3.Cells.Item(7, 7); $c.NumberFormat = "@"; $c.Value = '0.0161'
$q3.Cells.Item(7, 8).Value = 6
$ac = $q3.Cells.Item(8, 1)
$ac.Value = 6
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(8, 2); $c.NumberFormat = "@"; $c.Value = '009613'
$q3.Cells.Item(8, 3).Value = '上银中证500指数增强A'
$c = $q3.Cells.Item(8, 4); $c.NumberFormat = "@"; $c.Value = '0.97'
$c = $q3.Cells.Item(8, 5); $c.NumberFormat = "@"; $c.Value = '92.48'
$c = $q3.Cells.Item(8, 6); $c.NumberFormat = "@"; $c.Value = '1.00'
$c = $q3.Cells.Item(8, 7); $c.NumberFormat = "@"; $c.Value = '0.0097'
$q3.Cells.Item(8, 8).Value = 8
$ac = $q3.Cells.Item(9, 1)
$ac.Value = 7
$ac.Font.Bold = $true
$ac.HorizontalAlignment = -4108
$ac.VerticalAlignment = -4160
$c = $q3.Cells.Item(9, 2); $c.NumberFormat = "@"; $c.Value = '009614'
$q3.Cells.Item(9, 3).Value = '上银中证500指数增强C'
$c = $q3.Cells.Item(9, 4); $c.NumberFormat = "@"; $c.Value = '0.74'
$c = $q3.Cells.Item(9, 5); $c.NumberFormat = "@"; $c.Value = '92.48'
$c = $q3.Cells.Item(9, 6); $c.NumberFormat = "@"; $c.Value = '1.00'
$c = $q3.Cells.Item(9, 7); $c.NumberFormat = "@"; $c.Value = '0.0074'
$q3.Cells.Item(9, 8).Value = 8

# ------------------------------------------------------------------
# 2. Update the '总计' (grand total) sheet: insert a new row for
#    2022-Q3 right under the header and shift the existing quarters
#    down by one, keeping the running index in column A sequential.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$a = $total.Cells.Item(2, 1)
$a.Value = 0
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4160
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 0.24

# Renumber the running index (column A) of the rows that shifted down
for ($r = 3; $r -le 8; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ------------------------------------------------------------------
# 3. Keep '2020-Q4' as the selected/active sheet, matching the
#    original workbook's tab selection.
# ------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
